$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 442.5
$ws.Range("J28").Value = 429
$ws.Range("L28").Value = 429
$ws.Range("N28").Value = -1399
$ws.Range("H31").Value = 59824.332
$ws.Range("I31").Value = 71779.2
$ws.Range("K31").Value = 215337.6
$ws.Range("M31").Value = -215107.6
$ws.Range("H48").Value = 1500
$ws.Range("J48").Value = 1500
$ws.Range("L48").Value = 4500
$ws.Range("N48").Value = -5084
$ws.Range("H56").Value = 1500
$ws.Range("J56").Value = 1500
$ws.Range("L56").Value = 4500
$ws.Range("N56").Value = -5568
$ws.Range("H58").Value = 427.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 14998.75
$ws.Range("J86").Value = 14998.75
$ws.Range("L86").Value = 14998.75
$ws.Range("N86").Value = -17244.75
$ws.Range("H89").Value = 14998.75
$ws.Range("J89").Value = 14998.75
$ws.Range("L89").Value = 74993.75
$ws.Range("N89").Value = -86225.75
$ws.Range("H137").Value = 1035.8667
$ws.Range("I137").Value = 771.125
$ws.Range("K137").Value = 2313.375
$ws.Range("M137").Value = 236.625
$ws.Range("H138").Value = 4926.9653
$ws.Range("I138").Value = 3738.1667
$ws.Range("K138").Value = 11214.5001
$ws.Range("M138").Value = -6074.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H43").Value = 19974.5
$ws.Range("J43").Value = 19974.5
$ws.Range("L43").Value = 19974.5
$ws.Range("N43").Value = -20600.5
$ws.Range("H46").Value = 9561.25
$ws.Range("J46").Value = 9561.25
$ws.Range("L46").Value = 9561.25
$ws.Range("N46").Value = -10199.25
$ws.Range("H61").Value = 497.5
$ws.Range("I61").Value = 497.5
$ws.Range("K61").Value = 497.5
$ws.Range("M61").Value = -285.5
$ws.Range("H74").Value = 980.0909
$ws.Range("I74").Value = 1002.1667
$ws.Range("K74").Value = 1002.1667
$ws.Range("M74").Value = -128.1667
$ws.Range("H77").Value = 980.0909
$ws.Range("I77").Value = 1002.1667
$ws.Range("K77").Value = 5010.8335
$ws.Range("M77").Value = -642.8334999999997
$ws.Range("H88").Value = 2903.8
$ws.Range("I88").Value = 2837.75
$ws.Range("J88").Value = 2947.8333
$ws.Range("K88").Value = 2837.75
$ws.Range("L88").Value = 2947.8333
$ws.Range("M88").Value = -2431.75
$ws.Range("N88").Value = -3759.8333
$ws.Range("H91").Value = 2903.8
$ws.Range("I91").Value = 2837.75
$ws.Range("J91").Value = 2947.8333
$ws.Range("K91").Value = 2837.75
$ws.Range("L91").Value = 2947.8333
$ws.Range("M91").Value = -1433.75
$ws.Range("N91").Value = -5755.8333
$ws.Range("H97").Value = 660.9
$ws.Range("I97").Value = 551.25
$ws.Range("K97").Value = 551.25
$ws.Range("M97").Value = -55.25
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -470
$ws.Range("H136").Value = 497.5
$ws.Range("I136").Value = 497.5
$ws.Range("K136").Value = 1492.5
$ws.Range("M136").Value = 1057.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 333.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2221.6667
$ws.Range("J22").Value = 1733
$ws.Range("L22").Value = 1733
$ws.Range("N22").Value = -2433
$ws.Range("H31").Value = 1623.8108
$ws.Range("I31").Value = 1050.3846
$ws.Range("J31").Value = 1934.4166
$ws.Range("K31").Value = 1050.3846
$ws.Range("L31").Value = 1934.4166
$ws.Range("M31").Value = -755.3846000000001
$ws.Range("N31").Value = -2524.4166
$ws.Range("H34").Value = 1623.8108
$ws.Range("I34").Value = 1050.3846
$ws.Range("J34").Value = 1934.4166
$ws.Range("K34").Value = 1050.3846
$ws.Range("L34").Value = 1934.4166
$ws.Range("M34").Value = -848.3846000000001
$ws.Range("N34").Value = -2338.4166
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H99").Value = 1374.5555
$ws.Range("I99").Value = 1195
$ws.Range("J99").Value = 1599
$ws.Range("K99").Value = 1195
$ws.Range("L99").Value = 1599
$ws.Range("M99").Value = 303
$ws.Range("N99").Value = -4595
$ws.Range("H124").Value = 89717.8
$ws.Range("J124").Value = 89717.8
$ws.Range("L124").Value = 89717.8
$ws.Range("N124").Value = -94627.8
$ws.Range("H126").Value = 1374.5555
$ws.Range("I126").Value = 1195
$ws.Range("J126").Value = 1599
$ws.Range("K126").Value = 3585
$ws.Range("L126").Value = 4797
$ws.Range("M126").Value = -1115
$ws.Range("N126").Value = -9737
$ws.Range("H132").Value = 6826.75
$ws.Range("I132").Value = 6826.75
$ws.Range("K132").Value = 20480.25
$ws.Range("M132").Value = -17950.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 627.3333
$ws.Range("J34").Value = 832.6667
$ws.Range("L34").Value = 2498.0001
$ws.Range("N34").Value = -2666.0001
$ws.Range("H47").Value = 1027.3334
$ws.Range("I47").Value = 1027.3334
$ws.Range("K47").Value = 3082.0002
$ws.Range("M47").Value = -2651.0002
$ws.Range("H75").Value = 1458.5
$ws.Range("J75").Value = 419
$ws.Range("L75").Value = 1257
$ws.Range("N75").Value = -3253
$ws.Range("H78").Value = 1458.5
$ws.Range("J78").Value = 419
$ws.Range("L78").Value = 3771
$ws.Range("N78").Value = -13755
$ws.Range("H113").Value = 1065.5
$ws.Range("I113").Value = 384.25
$ws.Range("J113").Value = 1746.75
$ws.Range("K113").Value = 1152.75
$ws.Range("L113").Value = 5240.25
$ws.Range("M113").Value = 1017.25
$ws.Range("N113").Value = -9580.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33199.75
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 43766.332
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 43766.332
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -45762.332
$ws.Range("H83").Value = 33199.75
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 43766.332
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 218831.66
$ws.Range("M83").Value = -2508
$ws.Range("N83").Value = -228815.66
$ws.Range("H102").Value = 2403.75
$ws.Range("I102").Value = 2403.75
$ws.Range("K102").Value = 2403.75
$ws.Range("M102").Value = -781.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5906.933
$ws.Range("I7").Value = 3829.7144
$ws.Range("K7").Value = 3829.7144
$ws.Range("M7").Value = -3717.7144
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H46").Value = 1784.909
$ws.Range("I46").Value = 2146.5
$ws.Range("J46").Value = 1578.2858
$ws.Range("K46").Value = 2146.5
$ws.Range("L46").Value = 1578.2858
$ws.Range("M46").Value = -1958.5
$ws.Range("N46").Value = -1954.2858
$ws.Range("H100").Value = 2500
$ws.Range("I100").Value = 2500
$ws.Range("K100").Value = 2500
$ws.Range("M100").Value = -1959
$ws.Range("H126").Value = 5906.933
$ws.Range("I126").Value = 3829.7144
$ws.Range("K126").Value = 11489.1432
$ws.Range("M126").Value = -9019.143199999999
$ws.Range("H132").Value = 3039.3125
$ws.Range("I132").Value = 2941.889
$ws.Range("J132").Value = 3164.5715
$ws.Range("K132").Value = 8825.667000000001
$ws.Range("L132").Value = 9493.7145
$ws.Range("M132").Value = -6295.667000000001
$ws.Range("N132").Value = -14553.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1431050
$ws.Range("J81").Value = 2002900.4
$ws.Range("L81").Value = 4005800.8
$ws.Range("N81").Value = -4007922.8
$ws.Range("H84").Value = 1431050
$ws.Range("J84").Value = 2002900.4
$ws.Range("L84").Value = 20029004
$ws.Range("N84").Value = -20039612
$ws.Range("H122").Value = 4279.077
$ws.Range("I122").Value = 2673.111
$ws.Range("K122").Value = 8019.333
$ws.Range("M122").Value = -5569.333
